$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1746.25
$ws.Range("I86").Value = 1919.4445
$ws.Range("J86").Value = 1604.5454
$ws.Range("K86").Value = 1919.4445
$ws.Range("L86").Value = 1604.5454
$ws.Range("M86").Value = -796.4445000000001
$ws.Range("N86").Value = -3850.5454
$ws.Range("H89").Value = 1746.25
$ws.Range("I89").Value = 1919.4445
$ws.Range("J89").Value = 1604.5454
$ws.Range("K89").Value = 9597.2225
$ws.Range("L89").Value = 8022.727
$ws.Range("M89").Value = -3981.2225
$ws.Range("N89").Value = -19254.727
$ws.Range("H116").Value = 4478.75
$ws.Range("I116").Value = 3990
$ws.Range("J116").Value = 4772
$ws.Range("K116").Value = 3990
$ws.Range("L116").Value = 4772
$ws.Range("M116").Value = -548
$ws.Range("N116").Value = -11656
$ws.Range("H134").Value = 25954
$ws.Range("J134").Value = 25954
$ws.Range("L134").Value = 25954
$ws.Range("N134").Value = -36094
$ws.Range("H135").Value = 587
$ws.Range("I135").Value = 346.08694
$ws.Range("K135").Value = 3114.78246
$ws.Range("M135").Value = -579.7824600000004
$ws.Range("H137").Value = 3657
$ws.Range("I137").Value = 4165.25
$ws.Range("J137").Value = 2708.2666
$ws.Range("K137").Value = 12495.75
$ws.Range("L137").Value = 8124.7998
$ws.Range("M137").Value = -9945.75
$ws.Range("N137").Value = -13224.7998
$ws.Range("H141").Value = 1012435.9
$ws.Range("I141").Value = 1369.75
$ws.Range("J141").Value = 1686480
$ws.Range("K141").Value = 4109.25
$ws.Range("L141").Value = 5059440
$ws.Range("M141").Value = 1070.75
$ws.Range("N141").Value = -5069800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 28655.264
$ws.Range("I86").Value = 2010
$ws.Range("J86").Value = 128575
$ws.Range("K86").Value = 2010
$ws.Range("L86").Value = 128575
$ws.Range("M86").Value = -887
$ws.Range("N86").Value = -130821
$ws.Range("H89").Value = 28655.264
$ws.Range("I89").Value = 2010
$ws.Range("J89").Value = 128575
$ws.Range("K89").Value = 10050
$ws.Range("L89").Value = 642875
$ws.Range("M89").Value = -4434
$ws.Range("N89").Value = -654107

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1484.4166
$ws.Range("I16").Value = 620
$ws.Range("J16").Value = 2348.8333
$ws.Range("K16").Value = 620
$ws.Range("L16").Value = 2348.8333
$ws.Range("M16").Value = -333
$ws.Range("N16").Value = -2922.8333
$ws.Range("H22").Value = 890.0526
$ws.Range("I22").Value = 280.08334
$ws.Range("J22").Value = 1935.7142
$ws.Range("K22").Value = 280.08334
$ws.Range("L22").Value = 1935.7142
$ws.Range("M22").Value = 69.91665999999998
$ws.Range("N22").Value = -2635.7142
$ws.Range("H31").Value = 1778.59
$ws.Range("J31").Value = 2708.8696
$ws.Range("L31").Value = 2708.8696
$ws.Range("N31").Value = -3298.8696
$ws.Range("H34").Value = 1778.59
$ws.Range("J34").Value = 2708.8696
$ws.Range("L34").Value = 2708.8696
$ws.Range("N34").Value = -3112.8696
$ws.Range("H62").Value = 3941.8572
$ws.Range("I62").Value = 2450
$ws.Range("K62").Value = 2450
$ws.Range("M62").Value = -1826
$ws.Range("H65").Value = 3941.8572
$ws.Range("I65").Value = 2450
$ws.Range("K65").Value = 12250
$ws.Range("M65").Value = -9130
$ws.Range("H113").Value = 1484.4166
$ws.Range("I113").Value = 620
$ws.Range("J113").Value = 2348.8333
$ws.Range("K113").Value = 620
$ws.Range("L113").Value = 2348.8333
$ws.Range("M113").Value = 1550
$ws.Range("N113").Value = -6688.8333
$ws.Range("H122").Value = 2780.4783
$ws.Range("I122").Value = 2349.8667
$ws.Range("J122").Value = 3587.875
$ws.Range("K122").Value = 7049.6001
$ws.Range("L122").Value = 10763.625
$ws.Range("M122").Value = -4599.6001
$ws.Range("N122").Value = -15663.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 975.2222
$ws.Range("I18").Value = 445
$ws.Range("J18").Value = 1399.4
$ws.Range("K18").Value = 1335
$ws.Range("L18").Value = 4198.200000000001
$ws.Range("M18").Value = -1166
$ws.Range("N18").Value = -4536.200000000001
$ws.Range("H23").Value = 561.8182
$ws.Range("I23").Value = 89
$ws.Range("J23").Value = 666.8889
$ws.Range("K23").Value = 267
$ws.Range("L23").Value = 2000.6667
$ws.Range("M23").Value = -32
$ws.Range("N23").Value = -2470.6667
$ws.Range("H37").Value = 50850
$ws.Range("J37").Value = 50850
$ws.Range("L37").Value = 152550
$ws.Range("N37").Value = -152774
$ws.Range("H38").Value = 368.84616
$ws.Range("I38").Value = 74.27273
$ws.Range("J38").Value = 584.86664
$ws.Range("K38").Value = 222.81819
$ws.Range("L38").Value = 1754.59992
$ws.Range("M38").Value = 124.18181
$ws.Range("N38").Value = -2448.59992
$ws.Range("H87").Value = 12087.375
$ws.Range("I87").Value = 8224.75
$ws.Range("K87").Value = 24674.25
$ws.Range("M87").Value = -23426.25
$ws.Range("H90").Value = 12087.375
$ws.Range("I90").Value = 8224.75
$ws.Range("K90").Value = 74022.75
$ws.Range("M90").Value = -67782.75
$ws.Range("H92").Value = 1093.1666
$ws.Range("J92").Value = 1003.4737
$ws.Range("L92").Value = 3010.4211
$ws.Range("N92").Value = -5506.4211
$ws.Range("H131").Value = 839.80414
$ws.Range("I131").Value = 397.66666
$ws.Range("J131").Value = 885.0227
$ws.Range("K131").Value = 1192.99998
$ws.Range("L131").Value = 2655.0681
$ws.Range("M131").Value = 3847.00002
$ws.Range("N131").Value = -12735.0681

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5210.816
$ws.Range("I122").Value = 4750.65
$ws.Range("J122").Value = 5722.1113
$ws.Range("K122").Value = 14251.95
$ws.Range("L122").Value = 17166.3339
$ws.Range("M122").Value = -11801.95
$ws.Range("N122").Value = -22066.3339
$ws.Range("H126").Value = 457514.7
$ws.Range("I126").Value = 1834.25
$ws.Range("J126").Value = 1004331.2
$ws.Range("K126").Value = 5502.75
$ws.Range("L126").Value = 3012993.6
$ws.Range("M126").Value = -3032.75
$ws.Range("N126").Value = -3017933.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2346.375
$ws.Range("I7").Value = 1793.8182
$ws.Range("K7").Value = 1793.8182
$ws.Range("M7").Value = -1681.8182
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H126").Value = 2346.375
$ws.Range("I126").Value = 1793.8182
$ws.Range("K126").Value = 5381.4546
$ws.Range("M126").Value = -2911.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12138.659
$ws.Range("I132").Value = 3662.2646
$ws.Range("J132").Value = 34307.69
$ws.Range("K132").Value = 10986.7938
$ws.Range("L132").Value = 102923.07
$ws.Range("M132").Value = -8456.793799999999
$ws.Range("N132").Value = -107983.07
